$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to literal text while reusing an existing cell's style
# (so the style index matches exactly what Excel would produce, instead of
# letting a fresh NumberFormat assignment synthesize a brand-new style).
function Set-TextCell($ws, $addr, $text, $srcAddr) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($srcAddr).Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}

# --- Header / title updates -------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  44"
$ws.Range("C9").Value = "Report Covering the Week  10/30/2023  Through  11/5/2023"

# --- Weekly crime-statistics table updates (rows 14-29) --------------------
$ws.Range("N14").Value = -91.428571428571
$ws.Range("C15").Value = 1
$ws.Range("C15").NumberFormat = '#,##0'
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = '#,##0'
$ws.Range("E15").Value = 0
$ws.Range("E15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("F15").Value = 1
$ws.Range("F15").NumberFormat = '#,##0'
$ws.Range("G15").Value = 1
$ws.Range("G15").NumberFormat = '#,##0'
$ws.Range("H15").Value = 0
$ws.Range("H15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("I15").Value = 22
$ws.Range("J15").Value = 34
$ws.Range("K15").Value = -35.294117647058
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 4.761904761904
$ws.Range("N15").Value = -67.647058823529
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = -12.5
$ws.Range("I16").Value = 117
$ws.Range("J16").Value = 109
$ws.Range("K16").Value = 7.339449541284
$ws.Range("L16").Value = -7.142857142857
$ws.Range("M16").Value = -60.068259385665
$ws.Range("N16").Value = -87.671232876712
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 150
$ws.Range("F17").Value = 28
$ws.Range("G17").Value = 20
$ws.Range("H17").Value = 40
$ws.Range("I17").Value = 304
$ws.Range("J17").Value = 332
$ws.Range("K17").Value = -8.433734939759
$ws.Range("L17").Value = -10.059171597633
$ws.Range("M17").Value = 1.333333333333
$ws.Range("N17").Value = -52.351097178683
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -20
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 70
$ws.Range("J18").Value = 109
$ws.Range("K18").Value = -35.779816513761
$ws.Range("L18").Value = -36.936936936936
$ws.Range("M18").Value = -78.658536585365
$ws.Range("N18").Value = -92.600422832981
$ws.Range("C19").Value = 7
$ws.Range("C19").NumberFormat = '#,##0'
$ws.Range("E19").Value = -41.666666666666
$ws.Range("F19").Value = 17
$ws.Range("G19").Value = 44
$ws.Range("H19").Value = -61.363636363636
$ws.Range("I19").Value = 334
$ws.Range("J19").Value = 385
$ws.Range("K19").Value = -13.246753246753
$ws.Range("L19").Value = 14.383561643835
$ws.Range("M19").Value = -40.143369175627
$ws.Range("N19").Value = -90.290697674418
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = -66.666666666666
$ws.Range("F20").Value = 11
$ws.Range("G20").Value = 16
$ws.Range("H20").Value = -31.25
$ws.Range("I20").Value = 176
$ws.Range("J20").Value = 199
$ws.Range("K20").Value = -11.557788944723
$ws.Range("L20").Value = 18.120805369127
$ws.Range("M20").Value = -23.478260869565
$ws.Range("N20").Value = -88.179986568166
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = -7.692307692307
$ws.Range("F21").Value = 73
$ws.Range("G21").Value = 98
$ws.Range("H21").Value = -25.510204081632
$ws.Range("I21").Value = 1026
$ws.Range("J21").Value = 1176
$ws.Range("K21").Value = -12.755102040816
$ws.Range("L21").Value = -2.099236641221
$ws.Range("M21").Value = -41.270749856897
$ws.Range("N21").Value = -86.437541308658
Set-TextCell $ws "C23" "0" "D28"
$ws.Range("E23").Value = -100
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 100
$ws.Range("J23").Value = 6
$ws.Range("K23").Value = 133.333333333333
$ws.Range("M23").Value = 16.666666666666
$ws.Range("C24").Value = 25
$ws.Range("D24").Value = 26
$ws.Range("E24").Value = -3.846153846153
$ws.Range("F24").Value = 84
$ws.Range("G24").Value = 101
$ws.Range("H24").Value = -16.831683168316
$ws.Range("I24").Value = 985
$ws.Range("J24").Value = 1122
$ws.Range("K24").Value = -12.210338680926
$ws.Range("L24").Value = 22.057001239157
$ws.Range("M24").Value = 8.719646799117
$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = 60
$ws.Range("F25").Value = 58
$ws.Range("G25").Value = 36
$ws.Range("H25").Value = 61.111111111111
$ws.Range("I25").Value = 572
$ws.Range("J25").Value = 443
$ws.Range("K25").Value = 29.119638826185
$ws.Range("L25").Value = 36.842105263157
$ws.Range("M25").Value = -17.579250720461
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = -50
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = -50
$ws.Range("I26").Value = 39
$ws.Range("J26").Value = 48
$ws.Range("K26").Value = -18.75
$ws.Range("L26").Value = 0
Set-TextCell $ws "C27" "0" "D28"
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 8
$ws.Range("H27").Value = -75
$ws.Range("J27").Value = 51
$ws.Range("K27").Value = -21.568627450980
$ws.Range("L27").Value = -6.976744186046
$ws.Range("C28").Value = 1
$ws.Range("C28").NumberFormat = '#,##0'
$ws.Range("F28").Value = 2
Set-TextCell $ws "G28" "0" "D28"
Set-TextCell $ws "H28" "***.*" "E28"
$ws.Range("I28").Value = 18
$ws.Range("K28").Value = -40
$ws.Range("L28").Value = -60
$ws.Range("M28").Value = -64.705882352941
$ws.Range("N28").Value = -86.861313868613
$ws.Range("C29").Value = 1
$ws.Range("C29").NumberFormat = '#,##0'
$ws.Range("F29").Value = 2
Set-TextCell $ws "G29" "0" "D28"
Set-TextCell $ws "H29" "***.*" "E28"
$ws.Range("I29").Value = 13
$ws.Range("K29").Value = -45.833333333333
$ws.Range("L29").Value = -62.857142857142
$ws.Range("M29").Value = -69.047619047619
$ws.Range("N29").Value = -89.344262295082
